# Update saas_data.xlsx with newer data (adds July 2024 row to the sample-saas-data table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the month labels in column B (ene-24 -> enero, etc.) and add the
#    new "julio" label that will be used by the new row.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value2 = "enero"
$ws.Range("B3").Value2 = "febrero"
$ws.Range("B4").Value2 = "marzo"
$ws.Range("B5").Value2 = "abril"
$ws.Range("B6").Value2 = "mayo"
$ws.Range("B7").Value2 = "junio"

# ---------------------------------------------------------------------------
# 2. Add the new data row (row 8 = July 2024) to the table.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = 45474
$ws.Range("B8").Value2 = "julio"
$ws.Range("C8").Value2 = 4500000
$ws.Range("D8").Formula = "=(sample_saas_data[[#This Row],[arr]]-C7)/C7*100"
$ws.Range("E8").Value2 = 230
$ws.Range("F8").Formula = "=sample_saas_data[[#This Row],[arr]]/sample_saas_data[[#This Row],[customers]]"
$ws.Range("G8").Value2 = 50000
$ws.Range("H8").Value2 = 43000
$ws.Range("I8").Value2 = 9000
$ws.Range("J8").Value2 = -20000
$ws.Range("K8").Value2 = -9000

# Copy number formatting from the previous data row so row 8 looks the same
# as the rest of the table.
$ws.Range("A7:K7").Copy()
$ws.Range("A8:K8").PasteSpecial(-4122)

# Re-apply the values/formulas (PasteSpecial of formats only shouldn't have
# touched them, but make sure they are correct).
$ws.Range("A8").Value2 = 45474
$ws.Range("B8").Value2 = "julio"
$ws.Range("C8").Value2 = 4500000
$ws.Range("D8").Formula = "=(sample_saas_data[[#This Row],[arr]]-C7)/C7*100"
$ws.Range("E8").Value2 = 230
$ws.Range("F8").Formula = "=sample_saas_data[[#This Row],[arr]]/sample_saas_data[[#This Row],[customers]]"
$ws.Range("G8").Value2 = 50000
$ws.Range("H8").Value2 = 43000
$ws.Range("I8").Value2 = 9000
$ws.Range("J8").Value2 = -20000
$ws.Range("K8").Value2 = -9000

# ---------------------------------------------------------------------------
# 3. Resize the table / autofilter so it includes the new row.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K8"))

# ---------------------------------------------------------------------------
# 4. Add the formatted (still empty) row 9 below the table, mirroring the
#    formatting copied down from the data rows above.
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("F7:K7").Copy()
$ws.Range("F9:K9").PasteSpecial(-4122)

# D9 gets a percentage number format.
$ws.Range("D9").Style = "Porcentaje"
$ws.Range("D9").NumberFormat = "0.00%"

# B9 gets right alignment.
$ws.Range("B9").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 5. Tweak the number formats that differ between B2 (date-like format) and
#    B3:B8 (explicit General format).
# ---------------------------------------------------------------------------
$ws.Range("B2").NumberFormat = "mmm-yy"
$ws.Range("B3:B8").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 6. Update the defined name range used by the external data query.
# ---------------------------------------------------------------------------
$names = $wb.Names
$n = $names.Item(1)
$n.RefersTo = "='sample-saas-data'!`$A`$1:`$K`$8"

# ---------------------------------------------------------------------------
# 7. Update the active selection to match the author's final state.
# ---------------------------------------------------------------------------
$ws.Range("B9").Select()
